$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Append two new paragraphs after the document's final paragraph
# ("Climate Change Data from The World Bank: Data"):
#   1. an empty paragraph (spacing-after 0, no run)
#   2. a paragraph listing the countries (spacing-after 0)
# ---------------------------------------------------------------------

# Split off a new (still-empty) trailing paragraph after the current
# last paragraph.
$r1 = $d.Content
$r1.Collapse(0)
$r1.InsertParagraphAfter()

# Split off a second new (still-empty) trailing paragraph, so we end up
# with two blank paragraphs appended after the original last one.
$r2 = $d.Content
$r2.Collapse(0)
$r2.InsertParagraphAfter()

# Fill the very last (2nd new) paragraph with the countries text.
$r3 = $d.Content
$r3.Collapse(0)
$countries = "Argentina, Australia, Brazil, Canada, China, France, Germany, India, Indonesia, Italy, Japan, Mexico, Russian Federation, Saudi Arabia, South Africa, South Korea, Turkey, United Kingdom, United States, and the European Union"
$r3.InsertAfter($countries)

# `InsertParagraphAfter` leaves a stray empty <w:r/> in the first (blank)
# paragraph it created; clean that paragraph's XML so it matches a plain
# blank paragraph (spacing-after 0, no run at all).
$n = $d.Paragraphs.Count
$blankPara = $d.Paragraphs($n - 1)
$blankXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/></w:pPr></w:p>'
$blankPara.Range.InsertXML($blankXml)
